$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Date heading
Replace-Text "2024-06-19 Wednesday" "2024-06-20 Thursday"

# Table cell replacements, in document order.
# NOTE: "92×43=3956" -> "71×87=6177" and the original "71×87=6177" -> "11×92=1012"
# collide (the new text of the first equals the old text of the second), so the
# second one is applied first to avoid the Find matching the freshly written text.
Replace-Text "71×87=6177" "11×92=1012"

Replace-Text "87×65=5655" "69×62=4278"
Replace-Text "86×84=7224" "85×16=1360"
Replace-Text "37×40=1480" "34×96=3264"
Replace-Text "92×43=3956" "71×87=6177"
Replace-Text "84×78=6552" "72×79=5688"

Replace-Text "50×39=1950" "86×56=4816"
Replace-Text "68×89=6052" "64×16=1024"
Replace-Text "67×67=4489" "48×19=912"
Replace-Text "13×89=1157" "52×73=3796"
Replace-Text "96×99=9504" "70×45=3150"

Replace-Text "51×48=2448" "23×93=2139"
Replace-Text "73×64=4672" "93×27=2511"
Replace-Text "89×98=8722" "70×53=3710"
Replace-Text "53×97=5141" "97×95=9215"

Replace-Text "87×37=3219" "97×89=8633"
Replace-Text "29×18=522" "37×95=3515"
Replace-Text "82×69=5658" "97×73=7081"
Replace-Text "81×51=4131" "16×31=496"
Replace-Text "47×75=3525" "82×77=6314"

Replace-Text "38×19=722" "33×41=1353"
Replace-Text "64×21=1344" "29×15=435"
Replace-Text "84×79=6636" "45×15=675"
Replace-Text "76×35=2660" "62×47=2914"
Replace-Text "26×83=2158" "50×38=1900"
